$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update resistor rows (5-8) content: BOM re-grouping for v2.1 ---
# Row 5: Item labels for the 2.2K resistors change (R1,R2,R7,R8,R9 -> R1,R2,R4,R5,R6);
#        value/part info (2.2K) stays the same.
$ws.Range("C5").Value = "R1,R2,R4,R5,R6"

# Row 6: Item stays "R3" but its value changes from 10K to 1K (part info swaps in).
$ws.Range("D6").Value = "1K"
$ws.Range("E6").Value = "311-1.00KCRCT-ND"
$ws.Range("F6").Value = "RC0805FR-071KL"
$ws.Range("G6").Value = "RES 1K OHM 1% 1/8W 0805"

# Row 8: Item changes from "R5,R6" to "R9"; quantity 2 -> 1; value becomes 10K.
# (set before row 7 below so new shared strings land in the same order as the target file)
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "R9"
$ws.Range("D8").Value = "10K"
$ws.Range("E8").Value = "311-10.0KCRCT-ND"
$ws.Range("F8").Value = "RC0805FR-0710KL"
$ws.Range("G8").Value = "RES 10K OHM 1% 1/8W 0805"
$ws.Range("J8").Value = "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-0710KL/727535"

# Row 7: Item changes from "R4" to "R7,R8"; quantity 1 -> 2; value becomes 220 ohm.
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "R7,R8"
$ws.Range("D7").Value = 220
$ws.Range("E7").Value = "311-220ARCT-ND"
$ws.Range("F7").Value = "RC0805JR-07220RL"
$ws.Range("G7").Value = "RES 220 OHM 5% 1/8W 0805"
$ws.Range("J7").Value = "https://www.digikey.ca/en/products/detail/yageo/RC0805JR-07220RL/728283"

# Row 6: value's Digikey link text also needs to move to the 1K part page.
$ws.Range("J6").Value = "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-071KL/727444"

# --- Fix the total-price formula term order to match H8 before H6/H7 ---
$ws.Range("I13").Formula = "=H2*B2+H3*B3+H4*B4+H5*B5+H8*B8+H6*B6+H7*B7+H9*B9+H10*B10+H11*B11"

# --- Re-point the Digikey hyperlinks so each link follows its (possibly moved) part ---
$links = @(
  @{ Cell = "J2";  Url = "https://www.digikey.ca/en/products/detail/cui-devices/SJ1-3513/738683 " },
  @{ Cell = "J3";  Url = "https://www.digikey.ca/en/products/detail/texas-instruments/TCA9534DWR/6566100 " },
  @{ Cell = "J4";  Url = "https://www.digikey.ca/en/products/detail/sparkfun-electronics/PRT-14417/7652746 " },
  @{ Cell = "J5";  Url = "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-072K2L/727676" },
  @{ Cell = "J6";  Url = "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-071KL/727444" },
  @{ Cell = "J7";  Url = "https://www.digikey.ca/en/products/detail/yageo/RC0805JR-07220RL/728283" },
  @{ Cell = "J8";  Url = "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-0710KL/727535" },
  @{ Cell = "J9";  Url = "https://www.digikey.ca/en/products/detail/avx-corporation/08053C104KAT2A/1116281" },
  @{ Cell = "J10"; Url = "https://www.digikey.ca/en/products/detail/lite-on-inc/LTST-C190KRKT/386817" },
  @{ Cell = "J11"; Url = "https://www.digikey.ca/en/products/detail/omron-electronics-inc-emc-div/G3VM-41DY1-TR05/5799757 " }
)

$ws.Hyperlinks.Delete()
foreach ($link in $links) {
  $r = $ws.Range($link.Cell)
  $null = $ws.Hyperlinks.Add($r, $link.Url)
}

# --- Sheet view tidy-up: selection moved to E14, no fixed top-left cell anymore ---
$ws.Range("E14").Select()
